$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.372826457023621
$ws.Range("B1").Value = 2.645212173461914
$ws.Range("C1").Value = 5.745145320892334
$ws.Range("D1").Value = 2.235954761505127
$ws.Range("E1").Value = 1.219179153442383
